$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.505614041169197, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 2.210719231951476)
    3 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    4 = @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.613486173897131)
    5 = @(0.06328177979961902, 0.004309184025731883, 0.7127328510149897, 6.48142807727062, 7.261751892110961)
    6 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 11.47044854674929)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    8 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    9 = @(0.1554434735375247, 0.004309184025731883, 0.7127328510149897, 0.4998867070740569, 1.372372215652303)
    10 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    11 = @(0.3464964993005633, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.582219091977008)
    12 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    13 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    14 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    15 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.964442013611383)
    16 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 6.48142807727062, 12.7228780040422)
    17 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    18 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    19 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    20 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    21 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    22 = @(0.06328177979961902, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.328214039578707)
    23 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    24 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    25 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 20.64246832346449)
    26 = @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.613486173897131)
    27 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.705647867635037)
    28 = @(3.182878228561681, 0.05231270169004087, 0.7127328510149897, 6.48142807727062, 10.42935185853733)
    29 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    30 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    31 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    32 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    33 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    34 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    35 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    36 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    37 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    38 = @(0.7287194209349384, 1.65323645889881, 16.98373111632243, 6.48142807727062, 25.8471150734268)
    39 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
